# Text-to-speech V1 27 02 25
# Appends the new "text-to-speech" Q&A evaluation rows (tracks question,
# now including the "openai" model) to the filtered results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended after the existing 5 result rows (rows 2-6).
# Columns: A=index, B=Question, C=Model Name, D=Response, E=Expected Answer,
#          F=BLEU_score, G=ROUGE_score, H=similarity_score
$bleuTracks = [double]"2.570061783884003E-78"
$bleuCurves = [double]"6.373704167435469E-155"

$newRows = @(
    @{ A = 5; B = "How many tracks can you define in one ODF?"; C = "llama3.2:latest";
       D = "According to the Track Settings dialog box, the number of tracks that can be defined is 200.";
       E = "200 tracks can be defined in one ODF."; F = $bleuTracks; G = 0.32; H = 0.7143542766571045 },
    @{ A = 6; B = "How many tracks can you define in one ODF?"; C = "deepseek1.5";
       D = "According to the Track Settings dialog box, the number of tracks that can be defined is 200.";
       E = "200 tracks can be defined in one ODF."; F = $bleuTracks; G = 0.32; H = 0.7143542766571045 },
    @{ A = 7; B = "How many tracks can you define in one ODF?"; C = "openai";
       D = "According to the Track Settings dialog box, the number of tracks that can be defined is 200.";
       E = "200 tracks can be defined in one ODF."; F = $bleuTracks; G = 0.32; H = 0.7143542766571045 },
    @{ A = 8; B = "How many curves can I load in one go?"; C = "openai";
       D = "You can load up to 450 curves at a time.";
       E = "450 curves can be loaded in one go."; F = $bleuCurves; G = 0.2222222222222222; H = 0.9253911972045898 },
    @{ A = 9; B = "How many curves can I load in one go?"; C = "llama3.2:latest";
       D = "You can load up to 450 curves at a time.";
       E = "450 curves can be loaded in one go."; F = $bleuCurves; G = 0.2222222222222222; H = 0.9253911972045898 }
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H

    # Match the bordered/bold/centered style used on column A of the existing rows.
    $ws.Range("A6").Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
